$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode for row 4 (Forgot Password) from NO to YES (new distinct value "YES")
$ws.Range("C4").Value = "YES"

# Update the active selection to D8
$ws.Range("D8").Select()
